$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Collapse the split runs in the Title / Author / Abstract paragraphs into
#    single runs with the same (already-concatenated) text. Word's Find and
#    Replace naturally merges the destination into one run.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "Questions: Introduction to quadratic equations", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Questions: Introduction to quadratic equations", 2) | Out-Null

$d.Content.Find.Execute(
    "Tom Coleman", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Tom Coleman", 2) | Out-Null

$d.Content.Find.Execute(
    "A selection of questions for the study guide on introduction to quadratic equations.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A selection of questions for the study guide on introduction to quadratic equations.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Inside the equations, swap the order of <m:endChr/> and <m:sepChr/>
#    under every <m:dPr> (delimiter properties) so <m:sepChr/> precedes
#    <m:endChr/>. There is no OM property for this, so each affected
#    equation's OOXML is patched directly and written back with
#    Range.InsertXML, which replaces only that one equation's contents.
# ---------------------------------------------------------------------------

$begEndSep = '<m:begChr m:val="\(" ?/><m:endChr m:val="\)" ?/><m:sepChr m:val="" ?/><m:grow ?/>'
$begSepEnd = '<m:begChr m:val="("/><m:sepChr m:val=""/><m:endChr m:val=")"/><m:grow/>'

$count = $d.OMaths.Count
$prevStart = -1
$prevEnd = -1
$localIndex = 0

for ($i = 1; $i -le $count; $i++) {
    $curMath = $d.OMaths.Item($i)
    $curRange = $curMath.Range

    if ($curRange.Start -eq $prevStart -and $curRange.End -eq $prevEnd) {
        $localIndex = $localIndex + 1
    } else {
        $localIndex = 0
    }
    $prevStart = $curRange.Start
    $prevEnd = $curRange.End

    $packageXml = $curRange.WordOpenXML
    if ($packageXml -notlike "*m:sepChr*") {
        continue
    }

    $bodyStart = $packageXml.IndexOf("<w:body>")
    $bodyEnd = $packageXml.IndexOf("</w:body>")
    $body = $packageXml.Substring($bodyStart, ($bodyEnd - $bodyStart) + 9)

    $mathFragments = [regex]::Matches($body, '<m:oMath[^>]*>.*?</m:oMath>')
    if ($localIndex -ge $mathFragments.Count) {
        continue
    }

    $fragment = $mathFragments[$localIndex].Value
    if ($fragment -notlike "*m:sepChr*") {
        continue
    }

    $fixedFragment = [regex]::Replace($fragment, $begEndSep, $begSepEnd)
    if ($fixedFragment -ne $fragment) {
        $curRange.InsertXML($fixedFragment)
    }
}
